$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.639.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.124.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5276"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09096"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.125.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.851"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.109"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001182"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.012"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06732"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.334"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.702.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.388"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.361.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.575"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1082"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.664"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.410"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.025"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.941"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02660"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06898"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2327"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6935"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.276"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("E44").Value = "  +4.82%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6482"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.338"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000368"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.256"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07317"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
